# Apply the CV edits described by the commit:
# "Update CV with correct phone and SODA reviewing."

$d = $word.ActiveDocument

# 1. Phone number: split "+1 646 292 6641" into "+1 212 775" + " 6641"
#    (the diff splits it into two runs, but the visible text change is
#    simply replacing the old number with the new one).
$d.Content.Find.Execute("+1 646 292 6641", $false, $false, $false, $false,
                         $false, $true, 1, $false, "+1 212 775 6641", 2)

# 2. Add "SODA'17, " before "IEEE ICDM'16, " in the Conference Reviewing
#    list (use a right single quotation mark, matching the rest of the doc).
$soda = "SODA" + [char]0x2019 + "17, "
$icdm = "IEEE ICDM" + [char]0x2019 + "16, "
$d.Content.Find.Execute($icdm, $false, $false, $false, $false, $false,
                         $true, 1, $false, ($soda + $icdm), 2)

# 3. Update the footer page-count field's cached text from "1" to "7".
$footer = $d.Sections.First.Footers.First
$footer.Range.Find.Execute("1", $false, $false, $false, $false, $false,
                            $true, 1, $false, "7", 2)
